$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from an existing header cell (AC1) so the new header
# cells share the same bold/centered/bordered style used by the rest of
# row 1, then set the header text for the three new columns.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
$ws.Range("AD2:AD46").Value = 91
$ws.Range("AE2:AE46").Value = 71
$ws.Range("AF2:AF46").Value = 0
